# Rebuild Sheet1 content to match the updated restraints/demo export.
# The author's tool re-exported the sheet with:
#   - a new "box size" block (Ls=, Lx=, Ly=, Lz=) inserted above the table
#   - the "wj" column removed from the restraints table
#   - the restraints table shifted down to start at row 7
# We rebuild the sheet from scratch to avoid leaving stale cells behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all existing cell contents first so no stale cells (e.g. old row 6 / column H) remain.
$ws.Cells.ClearContents()

# --- Header ---
$ws.Range("A1").Value = "Restraints file"

# --- Box size summary block ---
$ws.Range("A2").Value = "Ls= "
$ws.Range("B2").Value = 10

$ws.Range("A3").Value = "Lx="
$ws.Range("B3").Value = 23.41

$ws.Range("A4").Value = "Ly="
$ws.Range("B4").Value = 18.690000000000001

$ws.Range("A5").Value = "Lz="
$ws.Range("B5").Value = 14.5

# Row 6 intentionally left blank.

# --- Restraints table header (row 7) ---
$ws.Range("A7").Value = "restraints"
$ws.Range("B7").Value = "prot x coor"
$ws.Range("C7").Value = "prot y coor"
$ws.Range("D7").Value = "prot z coor"
$ws.Range("E7").Value = "sl"
$ws.Range("F7").Value = "wi"
$ws.Range("G7").Value = "dij"

# --- Restraints table data (rows 8-10); "wj" column removed ---
$ws.Range("A8").Value = "36HE1-5HD2"
$ws.Range("B8").Value = 22.994400024413999
$ws.Range("C8").Value = 59.598503112792898
$ws.Range("D8").Value = 18.069942474365199
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 0.90909090909090906
$ws.Range("G8").Value = 2.899

$ws.Range("A9").Value = "54HH-7HB2"
$ws.Range("B9").Value = 29.420448303222599
$ws.Range("C9").Value = 66.502380371093693
$ws.Range("D9").Value = 22.527408599853501
$ws.Range("E9").Value = 7
$ws.Range("F9").Value = 0.81818181818181823
$ws.Range("G9").Value = 2.9289999999999998

$ws.Range("A10").Value = "8HE1-10HG2"
$ws.Range("B10").Value = 36.420028686523402
$ws.Range("C10").Value = 66.038200378417898
$ws.Range("D10").Value = 23.0035495758056
$ws.Range("E10").Value = 10
$ws.Range("F10").Value = 0.81818181818181823
$ws.Range("G10").Value = 2.1389999999999998

# Match the author's saved selection state.
$ws.Range("M10:M14").Select()
